# Re-computed NATMI TPM metrics for App-Cav1 ligand-receptor pairs.
# Ligand-side (G,H,I,J) values depend only on "Sending cluster" (col A);
# Receptor-side (M,N,O,P) values depend only on "Target cluster" (col D);
# K/L (receptor-expressing cells/rate) changed for Target cluster = Resolving-Mac.
# Edge weights: Q=G*M, R=H*N, S=I*O, T=J*P (recomputed automatically by Excel normally,
# but written explicitly here since this workbook stores static values, not formulas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs->ECs
$ws.Cells.Item(2, 7).Value = 89.38217433333334
$ws.Cells.Item(2, 8).Value = 268.146523
$ws.Cells.Item(2, 9).Value = 0.2143552015363441
$ws.Cells.Item(2, 10).Value = 0.2175965347165783
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 556.7425436666666
$ws.Cells.Item(2, 14).Value = 1670.227631
$ws.Cells.Item(2, 15).Value = 0.7235863858022448
$ws.Cells.Item(2, 16).Value = 0.7685368570853349
$ws.Cells.Item(2, 17).Value = 49762.85909679745
$ws.Cells.Item(2, 18).Value = 447865.731871177
$ws.Cells.Item(2, 19).Value = 0.155104505557595
$ws.Cells.Item(2, 20).Value = 0.1672309569037391

# Row 3: ECs->FAPs
$ws.Cells.Item(3, 7).Value = 89.38217433333334
$ws.Cells.Item(3, 8).Value = 268.146523
$ws.Cells.Item(3, 9).Value = 0.2143552015363441
$ws.Cells.Item(3, 10).Value = 0.2175965347165783
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 75.879851
$ws.Cells.Item(3, 14).Value = 227.639553
$ws.Cells.Item(3, 15).Value = 0.0986194207087145
$ws.Cells.Item(3, 16).Value = 0.1047458342586422
$ws.Cells.Item(3, 17).Value = 6782.306070469358
$ws.Cells.Item(3, 18).Value = 61040.75463422422
$ws.Cells.Item(3, 19).Value = 0.021139585801414
$ws.Cells.Item(3, 20).Value = 0.02279233056067759

# Row 4: ECs->Inflammatory-Mac
$ws.Cells.Item(4, 7).Value = 89.38217433333334
$ws.Cells.Item(4, 8).Value = 268.146523
$ws.Cells.Item(4, 9).Value = 0.2143552015363441
$ws.Cells.Item(4, 10).Value = 0.2175965347165783
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.402487333333333
$ws.Cells.Item(4, 14).Value = 4.207462
$ws.Cells.Item(4, 15).Value = 0.001822782814434402
$ws.Cells.Item(4, 16).Value = 0.001936017319896666
$ws.Cells.Item(4, 17).Value = 125.3573673282918
$ws.Cells.Item(4, 18).Value = 1128.216305954626
$ws.Cells.Item(4, 19).Value = 0.0003907229775450708
$ws.Cells.Item(4, 20).Value = 0.0004212706599607918

# Row 5: ECs->MuSCs
$ws.Cells.Item(5, 7).Value = 89.38217433333334
$ws.Cells.Item(5, 8).Value = 268.146523
$ws.Cells.Item(5, 9).Value = 0.2143552015363441
$ws.Cells.Item(5, 10).Value = 0.2175965347165783
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 135.006546
$ws.Cells.Item(5, 14).Value = 270.013092
$ws.Cells.Item(5, 15).Value = 0.1754651225976237
$ws.Cells.Item(5, 16).Value = 0.1242435517446983
$ws.Cells.Item(5, 17).Value = 12067.17863071319
$ws.Cells.Item(5, 18).Value = 72403.07178427912
$ws.Cells.Item(5, 19).Value = 0.03761186171701295
$ws.Cells.Item(5, 20).Value = 0.02703496632052625

# Row 6: ECs->Resolving-Mac
$ws.Cells.Item(6, 7).Value = 89.38217433333334
$ws.Cells.Item(6, 8).Value = 268.146523
$ws.Cells.Item(6, 9).Value = 0.2143552015363441
$ws.Cells.Item(6, 10).Value = 0.2175965347165783
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3895486666666667
$ws.Cells.Item(6, 14).Value = 1.168646
$ws.Cells.Item(6, 15).Value = 0.000506288076982634
$ws.Cells.Item(6, 16).Value = 0.0005377395914277917
$ws.Cells.Item(6, 17).Value = 34.81870683531756
$ws.Cells.Item(6, 18).Value = 313.368361517858
$ws.Cells.Item(6, 19).Value = 0.0001085254827770606
$ws.Cells.Item(6, 20).Value = 0.0001170102716745961

# Row 7: FAPs->ECs
$ws.Cells.Item(7, 7).Value = 122.3539896666667
$ws.Cells.Item(7, 8).Value = 367.061969
$ws.Cells.Item(7, 9).Value = 0.2934277926151677
$ws.Cells.Item(7, 10).Value = 0.2978648075949286
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 556.7425436666666
$ws.Cells.Item(7, 14).Value = 1670.227631
$ws.Cells.Item(7, 15).Value = 0.7235863858022448
$ws.Cells.Item(7, 16).Value = 0.7685368570853349
$ws.Cells.Item(7, 17).Value = 68119.67143478504
$ws.Cells.Item(7, 18).Value = 613077.0429130654
$ws.Cells.Item(7, 19).Value = 0.2123203559523398
$ws.Cells.Item(7, 20).Value = 0.2289200830653344

# Row 8: FAPs->FAPs
$ws.Cells.Item(8, 7).Value = 122.3539896666667
$ws.Cells.Item(8, 8).Value = 367.061969
$ws.Cells.Item(8, 9).Value = 0.2934277926151677
$ws.Cells.Item(8, 10).Value = 0.2978648075949286
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 75.879851
$ws.Cells.Item(8, 14).Value = 227.639553
$ws.Cells.Item(8, 15).Value = 0.0986194207087145
$ws.Cells.Item(8, 16).Value = 0.1047458342586422
$ws.Cells.Item(8, 17).Value = 9284.202505162206
$ws.Cells.Item(8, 18).Value = 83557.82254645985
$ws.Cells.Item(8, 19).Value = 0.02893767892754465
$ws.Cells.Item(8, 20).Value = 0.03120009776782072

# Row 9: FAPs->Inflammatory-Mac
$ws.Cells.Item(9, 7).Value = 122.3539896666667
$ws.Cells.Item(9, 8).Value = 367.061969
$ws.Cells.Item(9, 9).Value = 0.2934277926151677
$ws.Cells.Item(9, 10).Value = 0.2978648075949286
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.402487333333333
$ws.Cells.Item(9, 14).Value = 4.207462
$ws.Cells.Item(9, 15).Value = 0.001822782814434402
$ws.Cells.Item(9, 16).Value = 0.001936017319896666
$ws.Cells.Item(9, 17).Value = 171.5999206902975
$ws.Cells.Item(9, 18).Value = 1544.399286212678
$ws.Cells.Item(9, 19).Value = 0.0005348551376563494
$ws.Cells.Item(9, 20).Value = 0.0005766714264914696

# Row 10: FAPs->MuSCs
$ws.Cells.Item(10, 7).Value = 122.3539896666667
$ws.Cells.Item(10, 8).Value = 367.061969
$ws.Cells.Item(10, 9).Value = 0.2934277926151677
$ws.Cells.Item(10, 10).Value = 0.2978648075949286
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 135.006546
$ws.Cells.Item(10, 14).Value = 270.013092
$ws.Cells.Item(10, 15).Value = 0.1754651225976237
$ws.Cells.Item(10, 16).Value = 0.1242435517446983
$ws.Cells.Item(10, 17).Value = 16518.58953421636
$ws.Cells.Item(10, 18).Value = 99111.53720529815
$ws.Cells.Item(10, 19).Value = 0.05148634360477049
$ws.Cells.Item(10, 20).Value = 0.03700778163534512

# Row 11: FAPs->Resolving-Mac
$ws.Cells.Item(11, 7).Value = 122.3539896666667
$ws.Cells.Item(11, 8).Value = 367.061969
$ws.Cells.Item(11, 9).Value = 0.2934277926151677
$ws.Cells.Item(11, 10).Value = 0.2978648075949286
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3895486666666667
$ws.Cells.Item(11, 14).Value = 1.168646
$ws.Cells.Item(11, 15).Value = 0.000506288076982634
$ws.Cells.Item(11, 16).Value = 0.0005377395914277917
$ws.Cells.Item(11, 17).Value = 47.66283353599712
$ws.Cells.Item(11, 18).Value = 428.965501823974
$ws.Cells.Item(11, 19).Value = 0.0001485589928563923
$ws.Cells.Item(11, 20).Value = 0.0001601736999368147

# Row 12: Inflammatory-Mac->ECs
$ws.Cells.Item(12, 7).Value = 90.33462533333334
$ws.Cells.Item(12, 8).Value = 271.003876
$ws.Cells.Item(12, 9).Value = 0.2166393574945233
$ws.Cells.Item(12, 10).Value = 0.2199152301234996
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 556.7425436666666
$ws.Cells.Item(12, 14).Value = 1670.227631
$ws.Cells.Item(12, 15).Value = 0.7235863858022448
$ws.Cells.Item(12, 16).Value = 0.7685368570853349
$ws.Cells.Item(12, 17).Value = 50293.1290892553
$ws.Cells.Item(12, 18).Value = 452638.1618032978
$ws.Cells.Item(12, 19).Value = 0.1567572897119825
$ws.Cells.Item(12, 20).Value = 0.1690129597843126

# Row 13: Inflammatory-Mac->FAPs
$ws.Cells.Item(13, 7).Value = 90.33462533333334
$ws.Cells.Item(13, 8).Value = 271.003876
$ws.Cells.Item(13, 9).Value = 0.2166393574945233
$ws.Cells.Item(13, 10).Value = 0.2199152301234996
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 75.879851
$ws.Cells.Item(13, 14).Value = 227.639553
$ws.Cells.Item(13, 15).Value = 0.0986194207087145
$ws.Cells.Item(13, 16).Value = 0.1047458342586422
$ws.Cells.Item(13, 17).Value = 6854.577910434159
$ws.Cells.Item(13, 18).Value = 61691.20119390743
$ws.Cells.Item(13, 19).Value = 0.02136484793881799
$ws.Cells.Item(13, 20).Value = 0.02303520424546724

# Row 14: Inflammatory-Mac->Inflammatory-Mac
$ws.Cells.Item(14, 7).Value = 90.33462533333334
$ws.Cells.Item(14, 8).Value = 271.003876
$ws.Cells.Item(14, 9).Value = 0.2166393574945233
$ws.Cells.Item(14, 10).Value = 0.2199152301234996
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.402487333333333
$ws.Cells.Item(14, 14).Value = 4.207462
$ws.Cells.Item(14, 15).Value = 0.001822782814434402
$ws.Cells.Item(14, 16).Value = 0.001936017319896666
$ws.Cells.Item(14, 17).Value = 126.6931677914124
$ws.Cells.Item(14, 18).Value = 1140.238510122712
$ws.Cells.Item(14, 19).Value = 0.0003948864977711278
$ws.Cells.Item(14, 20).Value = 0.0004257596944281562

# Row 15: Inflammatory-Mac->MuSCs
$ws.Cells.Item(15, 7).Value = 90.33462533333334
$ws.Cells.Item(15, 8).Value = 271.003876
$ws.Cells.Item(15, 9).Value = 0.2166393574945233
$ws.Cells.Item(15, 10).Value = 0.2199152301234996
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 135.006546
$ws.Cells.Item(15, 14).Value = 270.013092
$ws.Cells.Item(15, 15).Value = 0.1754651225976237
$ws.Cells.Item(15, 16).Value = 0.1242435517446983
$ws.Cells.Item(15, 17).Value = 12195.76575045743
$ws.Cells.Item(15, 18).Value = 73174.5945027446
$ws.Cells.Item(15, 19).Value = 0.03801265142224695
$ws.Cells.Item(15, 20).Value = 0.02732304927329627

# Row 16: Inflammatory-Mac->Resolving-Mac
$ws.Cells.Item(16, 7).Value = 90.33462533333334
$ws.Cells.Item(16, 8).Value = 271.003876
$ws.Cells.Item(16, 9).Value = 0.2166393574945233
$ws.Cells.Item(16, 10).Value = 0.2199152301234996
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.3895486666666667
$ws.Cells.Item(16, 14).Value = 1.168646
$ws.Cells.Item(16, 15).Value = 0.000506288076982634
$ws.Cells.Item(16, 16).Value = 0.0005377395914277917
$ws.Cells.Item(16, 17).Value = 35.1897328524329
$ws.Cells.Item(16, 18).Value = 316.707595671896
$ws.Cells.Item(16, 19).Value = 0.0001096819237046556
$ws.Cells.Item(16, 20).Value = 0.0001182571259953595

# Row 17: MuSCs->ECs
$ws.Cells.Item(17, 7).Value = 18.634161
$ws.Cells.Item(17, 8).Value = 37.268322
$ws.Cells.Item(17, 9).Value = 0.0446882095496985
$ws.Cells.Item(17, 10).Value = 0.03024263611988591
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 556.7425436666666
$ws.Cells.Item(17, 14).Value = 1670.227631
$ws.Cells.Item(17, 15).Value = 0.7235863858022448
$ws.Cells.Item(17, 16).Value = 0.7685368570853349
$ws.Cells.Item(17, 17).Value = 10374.4301942342
$ws.Cells.Item(17, 18).Value = 62246.58116540518
$ws.Cells.Item(17, 19).Value = 0.0323357800360397
$ws.Cells.Item(17, 20).Value = 0.02324258051355255

# Row 18: MuSCs->FAPs
$ws.Cells.Item(18, 7).Value = 18.634161
$ws.Cells.Item(18, 8).Value = 37.268322
$ws.Cells.Item(18, 9).Value = 0.0446882095496985
$ws.Cells.Item(18, 10).Value = 0.03024263611988591
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 75.879851
$ws.Cells.Item(18, 14).Value = 227.639553
$ws.Cells.Item(18, 15).Value = 0.0986194207087145
$ws.Cells.Item(18, 16).Value = 0.1047458342586422
$ws.Cells.Item(18, 17).Value = 1413.957360190011
$ws.Cells.Item(18, 18).Value = 8483.744161140066
$ws.Cells.Item(18, 19).Value = 0.004407125338300909
$ws.Cells.Item(18, 20).Value = 0.003167790150557995

# Row 19: MuSCs->Inflammatory-Mac
$ws.Cells.Item(19, 7).Value = 18.634161
$ws.Cells.Item(19, 8).Value = 37.268322
$ws.Cells.Item(19, 9).Value = 0.0446882095496985
$ws.Cells.Item(19, 10).Value = 0.03024263611988591
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.402487333333333
$ws.Cells.Item(19, 14).Value = 4.207462
$ws.Cells.Item(19, 15).Value = 0.001822782814434402
$ws.Cells.Item(19, 16).Value = 0.001936017319896666
$ws.Cells.Item(19, 17).Value = 26.134174769794
$ws.Cells.Item(19, 18).Value = 156.805048618764
$ws.Cells.Item(19, 19).Value = 0.00008145690037503377
$ws.Cells.Item(19, 20).Value = 0.00005855026732743162

# Row 20: MuSCs->MuSCs
$ws.Cells.Item(20, 7).Value = 18.634161
$ws.Cells.Item(20, 8).Value = 37.268322
$ws.Cells.Item(20, 9).Value = 0.0446882095496985
$ws.Cells.Item(20, 10).Value = 0.03024263611988591
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 135.006546
$ws.Cells.Item(20, 14).Value = 270.013092
$ws.Cells.Item(20, 15).Value = 0.1754651225976237
$ws.Cells.Item(20, 16).Value = 0.1242435517446983
$ws.Cells.Item(20, 17).Value = 2515.733714217906
$ws.Cells.Item(20, 18).Value = 10062.93485687163
$ws.Cells.Item(20, 19).Value = 0.007841222167306145
$ws.Cells.Item(20, 20).Value = 0.003757452525657129

# Row 21: MuSCs->Resolving-Mac
$ws.Cells.Item(21, 7).Value = 18.634161
$ws.Cells.Item(21, 8).Value = 37.268322
$ws.Cells.Item(21, 9).Value = 0.0446882095496985
$ws.Cells.Item(21, 10).Value = 0.03024263611988591
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.3895486666666667
$ws.Cells.Item(21, 14).Value = 1.168646
$ws.Cells.Item(21, 15).Value = 0.000506288076982634
$ws.Cells.Item(21, 16).Value = 0.0005377395914277917
$ws.Cells.Item(21, 17).Value = 7.258912572002
$ws.Cells.Item(21, 18).Value = 43.553475432012
$ws.Cells.Item(21, 19).Value = 0.00002262510767671383
$ws.Cells.Item(21, 20).Value = 0.00001626266279080683

# Row 22: Resolving-Mac->ECs
$ws.Cells.Item(22, 7).Value = 96.27664699999998
$ws.Cells.Item(22, 8).Value = 288.829941
$ws.Cells.Item(22, 9).Value = 0.2308894388042666
$ws.Cells.Item(22, 10).Value = 0.2343807914451077
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 556.7425436666666
$ws.Cells.Item(22, 14).Value = 1670.227631
$ws.Cells.Item(22, 15).Value = 0.7235863858022448
$ws.Cells.Item(22, 16).Value = 0.7685368570853349
$ws.Cells.Item(22, 17).Value = 53601.30534647774
$ws.Cells.Item(22, 18).Value = 482411.7481182997
$ws.Cells.Item(22, 19).Value = 0.1670684545442878
$ws.Cells.Item(22, 20).Value = 0.1801302768183964

# Row 23: Resolving-Mac->FAPs
$ws.Cells.Item(23, 7).Value = 96.27664699999998
$ws.Cells.Item(23, 8).Value = 288.829941
$ws.Cells.Item(23, 9).Value = 0.2308894388042666
$ws.Cells.Item(23, 10).Value = 0.2343807914451077
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 75.879851
$ws.Cells.Item(23, 14).Value = 227.639553
$ws.Cells.Item(23, 15).Value = 0.0986194207087145
$ws.Cells.Item(23, 16).Value = 0.1047458342586422
$ws.Cells.Item(23, 17).Value = 7305.457629139596
$ws.Cells.Item(23, 18).Value = 65749.11866225637
$ws.Cells.Item(23, 19).Value = 0.02277018270263696
$ws.Cells.Item(23, 20).Value = 0.02455041153411862

# Row 24: Resolving-Mac->Inflammatory-Mac
$ws.Cells.Item(24, 7).Value = 96.27664699999998
$ws.Cells.Item(24, 8).Value = 288.829941
$ws.Cells.Item(24, 9).Value = 0.2308894388042666
$ws.Cells.Item(24, 10).Value = 0.2343807914451077
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 1.402487333333333
$ws.Cells.Item(24, 14).Value = 4.207462
$ws.Cells.Item(24, 15).Value = 0.001822782814434402
$ws.Cells.Item(24, 16).Value = 0.001936017319896666
$ws.Cells.Item(24, 17).Value = 135.0267779133046
$ws.Cells.Item(24, 18).Value = 1215.241001219742
$ws.Cells.Item(24, 19).Value = 0.0004208613010868207
$ws.Cells.Item(24, 20).Value = 0.0004537652716888167

# Row 25: Resolving-Mac->MuSCs
$ws.Cells.Item(25, 7).Value = 96.27664699999998
$ws.Cells.Item(25, 8).Value = 288.829941
$ws.Cells.Item(25, 9).Value = 0.2308894388042666
$ws.Cells.Item(25, 10).Value = 0.2343807914451077
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 135.006546
$ws.Cells.Item(25, 14).Value = 270.013092
$ws.Cells.Item(25, 15).Value = 0.1754651225976237
$ws.Cells.Item(25, 16).Value = 0.1242435517446983
$ws.Cells.Item(25, 17).Value = 12997.97757193126
$ws.Cells.Item(25, 18).Value = 77987.86543158758
$ws.Cells.Item(25, 19).Value = 0.04051304368628716
$ws.Cells.Item(25, 20).Value = 0.02912030198987358

# Row 26: Resolving-Mac->Resolving-Mac
$ws.Cells.Item(26, 7).Value = 96.27664699999998
$ws.Cells.Item(26, 8).Value = 288.829941
$ws.Cells.Item(26, 9).Value = 0.2308894388042666
$ws.Cells.Item(26, 10).Value = 0.2343807914451077
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.3895486666666667
$ws.Cells.Item(26, 14).Value = 1.168646
$ws.Cells.Item(26, 15).Value = 0.000506288076982634
$ws.Cells.Item(26, 16).Value = 0.0005377395914277917
$ws.Cells.Item(26, 17).Value = 37.50443946998733
$ws.Cells.Item(26, 18).Value = 337.539955229886
$ws.Cells.Item(26, 19).Value = 0.0001168965699678117
$ws.Cells.Item(26, 20).Value = 0.0001260358310302146
